$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed in B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values changed in B2:E2
$ws.Range("B2").Value = 35.858458423859474
$ws.Range("C2").Value = 43.708302383168899
$ws.Range("D2").Value = 39.514648813360758
$ws.Range("E2").Value = 45.014378046495629

# Row 3 values changed in B3:E3
$ws.Range("B3").Value = 42.42826502455631
$ws.Range("C3").Value = 51.213714597804675
$ws.Range("D3").Value = 49.837777848804549
$ws.Range("E3").Value = 44.238443142811157

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
